$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2022-07-11 20:58:10"

for ($row = 2; $row -le 35; $row++) {
    $ws.Cells.Item($row, 15).Value = $newTimestamp
}
